$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''28.156.04'
$ws.Range("E2").Value = '  +0.04%  '

$ws.Range("D3").Value = '''1.826.50'
$ws.Range("E3").Value = '  +1.87%  '

$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  -0.69%  '

$ws.Range("E5").Value = '  -0.52%  '

$ws.Range("E6").Value = '  -0.53%  '

$ws.Range("D7").Value = '''0.5131'
$ws.Range("E7").Value = '  -2.00%  '

$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D8").Value = '''0.1009'
$ws.Range("E8").Value = '  +26.95%  '

$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '''0.3925'
$ws.Range("E9").Value = '  +3.10%  '

$ws.Range("D10").Value = '''1.112'
$ws.Range("E10").Value = '  +1.78%  '

$ws.Range("D11").Value = '''41.00'
$ws.Range("E11").Value = '  -0.85%  '

$ws.Range("D12").Value = '''6.490'
$ws.Range("E12").Value = '  +3.81%  '

$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").Value = '''20.84'
$ws.Range("E13").Value = '  +1.65%  '

$ws.Range("B14").Value = 'BinanceUSD'
$ws.Range("C14").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D14").Value = '''1.0000'
$ws.Range("E14").Value = '  -0.58%  '

$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '''1.833.27'
$ws.Range("E15").Value = '  +1.77%  '

$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '''7.410'
$ws.Range("E16").Value = '  +1.48%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '''0.00001142'
$ws.Range("E17").Value = '  +5.04%  '

$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").Value = '''94.53'
$ws.Range("E18").Value = '  +3.22%  '

$ws.Range("D19").Value = '''0.06612'
$ws.Range("E19").Value = '  +0.37%  '

$ws.Range("E20").Value = '  -0.47%  '

$ws.Range("D21").Value = '''17.37'
$ws.Range("E21").Value = '  +0.53%  '

$ws.Range("D22").Value = '''6.054'
$ws.Range("E22").Value = '  +1.68%  '

$ws.Range("D23").Value = '''28.229.95'
$ws.Range("E23").Value = '  +0.10%  '

$ws.Range("D24").Value = '''11.20'
$ws.Range("E24").Value = '  +0.88%  '

$ws.Range("D25").Value = '''2.247'
$ws.Range("E25").Value = '  -0.94%  '

$ws.Range("D26").Value = '''159.01'
$ws.Range("E26").Value = '  -0.98%  '

$ws.Range("D27").Value = '''2.466'
$ws.Range("E27").Value = '  +5.58%  '

$ws.Range("D28").Value = '''20.86'
$ws.Range("E28").Value = '  +2.42%  '

$ws.Range("D29").Value = '''2.037.59'
$ws.Range("E29").Value = '  +1.67%  '

$ws.Range("D30").Value = '''128.97'
$ws.Range("E30").Value = '  +5.34%  '

$ws.Range("D31").Value = '''0.1094'
$ws.Range("E31").Value = '  +1.14%  '

$ws.Range("D32").Value = '''1.074'
$ws.Range("E32").Value = '  +2.22%  '

$ws.Range("D33").Value = '''5.651'
$ws.Range("E33").Value = '  +2.58%  '

$ws.Range("D34").Value = '''3.635'
$ws.Range("E34").Value = '  -1.53%  '

$ws.Range("D35").Value = '''0.06951'
$ws.Range("E35").Value = '  -3.75%  '

$ws.Range("D36").Value = '''9.164'
$ws.Range("E36").Value = '  +6.95%  '

$ws.Range("D37").Value = '''0.02351'
$ws.Range("E37").Value = '  +2.05%  '

$ws.Range("D38").Value = '''0.2179'
$ws.Range("E38").Value = '  +1.68%  '

$ws.Range("D39").Value = '''11.67'
$ws.Range("E39").Value = '  -4.27%  '

$ws.Range("E40").Value = '  -0.56%  '

$ws.Range("D41").Value = '''0.6285'
$ws.Range("E41").Value = '  +2.10%  '

$ws.Range("D42").Value = '''0.9998'
$ws.Range("E42").Value = '  -0.42%  '

$ws.Range("D43").Value = '''1.160'
$ws.Range("E43").Value = '  -0.25%  '

$ws.Range("D44").Value = '''13.38'
$ws.Range("E44").Value = '  +0.78%  '

$ws.Range("D45").Value = '''0.6010'
$ws.Range("E45").Value = '  +0.78%  '

$ws.Range("D46").Value = '''3.724'
$ws.Range("E46").Value = '  -1.19%  '

$ws.Range("E47").Value = '  -6.23%  '

$ws.Range("D48").Value = '''125.96'
$ws.Range("E48").Value = '  -0.51%  '

$ws.Range("D49").Value = '''2.006'
$ws.Range("E49").Value = '  +4.75%  '

$ws.Range("E50").Value = '  -1.92%  '

$ws.Range("D51").Value = '''0.06787'
$ws.Range("E51").Value = '  +0.06%  '
